# Auto-generated script: refresh market-price derived columns (H-N)
# across all Leve-profit sheets, per the scheduled-runner data update.
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 3251.5588  # H132: 3328.879 -> 3251.5588
$ws.Cells.Item(132, 9).Value = 1516.0646  # I132: 1543.2667 -> 1516.0646
$ws.Cells.Item(132, 11).Value = 4548.1938  # K132: 4629.800099999999 -> 4548.1938
$ws.Cells.Item(132, 13).Value = -2018.1938  # M132: -2099.800099999999 -> -2018.1938
$ws.Cells.Item(136, 8).Value = 0  # H136: 150000 -> 0
$ws.Cells.Item(136, 10).Value = 0  # J136: 150000 -> 0
$ws.Cells.Item(136, 14).ClearContents()  # N136: remove cell (was -160200)
$ws.Cells.Item(138, 8).Value = 2366.4343  # H138: 2325.0845 -> 2366.4343
$ws.Cells.Item(138, 10).Value = 3203.5454  # J138: 3235.5898 -> 3203.5454
$ws.Cells.Item(138, 12).Value = 9610.636200000001  # L138: 9706.769400000001 -> 9610.636200000001
$ws.Cells.Item(138, 14).Value = -19890.6362  # N138: -19986.7694 -> -19890.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2815415  # H32: 2548694.2 -> 2815415
$ws.Cells.Item(32, 9).Value = 547966  # I32: 489258.34 -> 547966
$ws.Cells.Item(32, 11).Value = 547966  # K32: 489258.34 -> 547966
$ws.Cells.Item(32, 13).Value = -547679  # M32: -488971.34 -> -547679
$ws.Cells.Item(133, 8).Value = 0  # H133: 50000 -> 0
$ws.Cells.Item(133, 10).Value = 0  # J133: 50000 -> 0
$ws.Cells.Item(133, 14).ClearContents()  # N133: remove cell (was -55060)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(36, 8).Value = 6320  # H36: 8047.4546 -> 6320
$ws.Cells.Item(36, 9).Value = 1149.6666  # I36: 499.83334 -> 1149.6666
$ws.Cells.Item(36, 10).Value = 16660.666  # J36: 17104.6 -> 16660.666
$ws.Cells.Item(36, 11).Value = 1149.6666  # K36: 499.83334 -> 1149.6666
$ws.Cells.Item(36, 12).Value = 16660.666  # L36: 17104.6 -> 16660.666
$ws.Cells.Item(36, 13).Value = -615.6666  # M36: 34.16665999999998 -> -615.6666
$ws.Cells.Item(36, 14).Value = -17728.666  # N36: -18172.6 -> -17728.666
$ws.Cells.Item(86, 8).Value = 2996.3333  # H86: 2997.5 -> 2996.3333
$ws.Cells.Item(86, 9).Value = 2994  # I86: 0 -> 2994
$ws.Cells.Item(86, 11).Value = 2994  # K86: 0 -> 2994
$ws.Cells.Item(86, 13).Value = -1871  # M86: None -> -1871
$ws.Cells.Item(89, 8).Value = 2996.3333  # H89: 2997.5 -> 2996.3333
$ws.Cells.Item(89, 9).Value = 2994  # I89: 0 -> 2994
$ws.Cells.Item(89, 11).Value = 14970  # K89: 0 -> 14970
$ws.Cells.Item(89, 13).Value = -9354  # M89: None -> -9354
$ws.Cells.Item(105, 8).Value = 19121.24  # H105: 12529.795 -> 19121.24
$ws.Cells.Item(105, 9).Value = 2377.353  # I105: 1668.2 -> 2377.353
$ws.Cells.Item(105, 10).Value = 54702  # J105: 48735.11 -> 54702
$ws.Cells.Item(105, 11).Value = 2377.353  # K105: 1668.2 -> 2377.353
$ws.Cells.Item(105, 12).Value = 54702  # L105: 48735.11 -> 54702
$ws.Cells.Item(105, 13).Value = -630.3530000000001  # M105: 78.79999999999995 -> -630.3530000000001
$ws.Cells.Item(105, 14).Value = -58196  # N105: -52229.11 -> -58196
$ws.Cells.Item(134, 8).Value = 3289.923  # H134: 3447.9167 -> 3289.923
$ws.Cells.Item(134, 9).Value = 3317.147  # I134: 3356.6177 -> 3317.147
$ws.Cells.Item(134, 10).Value = 3104.8  # J134: 5000 -> 3104.8
$ws.Cells.Item(134, 11).Value = 9951.440999999999  # K134: 10069.8531 -> 9951.440999999999
$ws.Cells.Item(134, 12).Value = 9314.400000000001  # L134: 15000 -> 9314.400000000001
$ws.Cells.Item(134, 13).Value = -7416.440999999999  # M134: -7534.8531 -> -7416.440999999999
$ws.Cells.Item(134, 14).Value = -14384.4  # N134: -20070 -> -14384.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2543.2954  # H31: 2406.4666 -> 2543.2954
$ws.Cells.Item(31, 9).Value = 1408.1482  # I31: 1358.8276 -> 1408.1482
$ws.Cells.Item(31, 10).Value = 4346.1763  # J31: 4305.3125 -> 4346.1763
$ws.Cells.Item(31, 11).Value = 1408.1482  # K31: 1358.8276 -> 1408.1482
$ws.Cells.Item(31, 12).Value = 4346.1763  # L31: 4305.3125 -> 4346.1763
$ws.Cells.Item(31, 13).Value = -1113.1482  # M31: -1063.8276 -> -1113.1482
$ws.Cells.Item(31, 14).Value = -4936.1763  # N31: -4895.3125 -> -4936.1763
$ws.Cells.Item(34, 8).Value = 2543.2954  # H34: 2406.4666 -> 2543.2954
$ws.Cells.Item(34, 9).Value = 1408.1482  # I34: 1358.8276 -> 1408.1482
$ws.Cells.Item(34, 10).Value = 4346.1763  # J34: 4305.3125 -> 4346.1763
$ws.Cells.Item(34, 11).Value = 1408.1482  # K34: 1358.8276 -> 1408.1482
$ws.Cells.Item(34, 12).Value = 4346.1763  # L34: 4305.3125 -> 4346.1763
$ws.Cells.Item(34, 13).Value = -1206.1482  # M34: -1156.8276 -> -1206.1482
$ws.Cells.Item(34, 14).Value = -4750.1763  # N34: -4709.3125 -> -4750.1763
$ws.Cells.Item(141, 8).Value = 86666.5  # H141: 90999.836 -> 86666.5
$ws.Cells.Item(141, 10).Value = 86666.5  # J141: 90999.836 -> 86666.5
$ws.Cells.Item(141, 12).Value = 86666.5  # L141: 90999.836 -> 86666.5
$ws.Cells.Item(141, 14).Value = -97026.5  # N141: -101359.836 -> -97026.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 7321.05  # H69: 8279.210999999999 -> 7321.05
$ws.Cells.Item(69, 9).Value = 3202.4285  # I69: 3660.4 -> 3202.4285
$ws.Cells.Item(69, 10).Value = 9538.77  # J69: 9928.786 -> 9538.77
$ws.Cells.Item(69, 11).Value = 9607.2855  # K69: 10981.2 -> 9607.2855
$ws.Cells.Item(69, 12).Value = 28616.31  # L69: 29786.358 -> 28616.31
$ws.Cells.Item(69, 13).Value = -8796.2855  # M69: -10170.2 -> -8796.2855
$ws.Cells.Item(69, 14).Value = -30238.31  # N69: -31408.358 -> -30238.31
$ws.Cells.Item(72, 8).Value = 7321.05  # H72: 8279.210999999999 -> 7321.05
$ws.Cells.Item(72, 9).Value = 3202.4285  # I72: 3660.4 -> 3202.4285
$ws.Cells.Item(72, 10).Value = 9538.77  # J72: 9928.786 -> 9538.77
$ws.Cells.Item(72, 11).Value = 28821.8565  # K72: 32943.6 -> 28821.8565
$ws.Cells.Item(72, 12).Value = 85848.93000000001  # L72: 89359.07399999999 -> 85848.93000000001
$ws.Cells.Item(72, 13).Value = -24765.8565  # M72: -28887.6 -> -24765.8565
$ws.Cells.Item(72, 14).Value = -93960.93000000001  # N72: -97471.07399999999 -> -93960.93000000001
$ws.Cells.Item(99, 8).Value = 23996  # H99: 24000 -> 23996
$ws.Cells.Item(99, 10).Value = 23996  # J99: 24000 -> 23996
$ws.Cells.Item(99, 12).Value = 71988  # L99: 72000 -> 71988
$ws.Cells.Item(99, 14).Value = -76480  # N99: -76492 -> -76480
$ws.Cells.Item(102, 8).Value = 8811.909  # H102: 6713.222 -> 8811.909
$ws.Cells.Item(102, 9).Value = 5000  # I102: 0 -> 5000
$ws.Cells.Item(102, 10).Value = 9193.1  # J102: 6713.222 -> 9193.1
$ws.Cells.Item(102, 11).Value = 15000  # K102: 0 -> 15000
$ws.Cells.Item(102, 12).Value = 27579.3  # L102: 20139.666 -> 27579.3
$ws.Cells.Item(102, 13).Value = -12566  # M102: None -> -12566
$ws.Cells.Item(102, 14).Value = -32447.3  # N102: -25007.666 -> -32447.3
$ws.Cells.Item(104, 8).Value = 174037.33  # H104: 501499.5 -> 174037.33
$ws.Cells.Item(104, 9).Value = 335675  # I104: 501499.5 -> 335675
$ws.Cells.Item(104, 10).Value = 12399.667  # J104: 0 -> 12399.667
$ws.Cells.Item(104, 11).Value = 1007025  # K104: 1504498.5 -> 1007025
$ws.Cells.Item(104, 12).Value = 37199.001  # L104: 0 -> 37199.001
$ws.Cells.Item(104, 13).Value = -1004404  # M104: -1501877.5 -> -1004404
$ws.Cells.Item(104, 14).Value = -42441.001  # N104: None -> -42441.001
$ws.Cells.Item(107, 8).Value = 729.7143  # H107: 734.8 -> 729.7143
$ws.Cells.Item(107, 10).Value = 728.2353000000001  # J107: 734.5 -> 728.2353000000001
$ws.Cells.Item(107, 12).Value = 2184.7059  # L107: 2203.5 -> 2184.7059
$ws.Cells.Item(107, 14).Value = -6024.7059  # N107: -6043.5 -> -6024.7059
$ws.Cells.Item(134, 8).Value = 1311.25  # H134: 1551 -> 1311.25
$ws.Cells.Item(134, 9).Value = 1370  # I134: 1551 -> 1370
$ws.Cells.Item(134, 10).Value = 900  # J134: 0 -> 900
$ws.Cells.Item(134, 11).Value = 4110  # K134: 4653 -> 4110
$ws.Cells.Item(134, 12).Value = 2700  # L134: 0 -> 2700
$ws.Cells.Item(134, 13).Value = 960  # M134: 417 -> 960
$ws.Cells.Item(134, 14).Value = -12840  # N134: None -> -12840
$ws.Cells.Item(140, 8).Value = 4696.684  # H140: 2937.6572 -> 4696.684
$ws.Cells.Item(140, 9).Value = 5212  # I140: 1792.762 -> 5212
$ws.Cells.Item(140, 10).Value = 4512.643  # J140: 4655 -> 4512.643
$ws.Cells.Item(140, 11).Value = 15636  # K140: 5378.286 -> 15636
$ws.Cells.Item(140, 12).Value = 13537.929  # L140: 13965 -> 13537.929
$ws.Cells.Item(140, 13).Value = -10456  # M140: -198.2860000000001 -> -10456
$ws.Cells.Item(140, 14).Value = -23897.929  # N140: -24325 -> -23897.929
$ws.Cells.Item(141, 8).Value = 7166.5  # H141: 4999.4 -> 7166.5
$ws.Cells.Item(141, 9).Value = 2999.6667  # I141: 2284.8572 -> 2999.6667
$ws.Cells.Item(141, 11).Value = 8999.000100000001  # K141: 6854.571599999999 -> 8999.000100000001
$ws.Cells.Item(141, 13).Value = -3819.000100000001  # M141: -1674.571599999999 -> -3819.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 25013324  # H107: 23822258 -> 25013324
$ws.Cells.Item(107, 9).Value = 1154.909  # I107: 1137.75 -> 1154.909
$ws.Cells.Item(107, 11).Value = 1154.909  # K107: 1137.75 -> 1154.909
$ws.Cells.Item(107, 13).Value = 765.0909999999999  # M107: 782.25 -> 765.0909999999999
$ws.Cells.Item(113, 8).Value = 11210.954  # H113: 10841 -> 11210.954
$ws.Cells.Item(113, 9).Value = 14304.5  # I113: 14379.625 -> 14304.5
$ws.Cells.Item(113, 10).Value = 2961.5  # J113: 2752.7144 -> 2961.5
$ws.Cells.Item(113, 11).Value = 14304.5  # K113: 14379.625 -> 14304.5
$ws.Cells.Item(113, 12).Value = 2961.5  # L113: 2752.7144 -> 2961.5
$ws.Cells.Item(113, 13).Value = -12134.5  # M113: -12209.625 -> -12134.5
$ws.Cells.Item(113, 14).Value = -7301.5  # N113: -7092.7144 -> -7301.5
$ws.Cells.Item(122, 8).Value = 1568.2  # H122: 1608 -> 1568.2
$ws.Cells.Item(122, 9).Value = 1458.25  # I122: 1474.6666 -> 1458.25
$ws.Cells.Item(122, 11).Value = 4374.75  # K122: 4423.9998 -> 4374.75
$ws.Cells.Item(122, 13).Value = -1924.75  # M122: -1973.9998 -> -1924.75
$ws.Cells.Item(126, 8).Value = 44610.734  # H126: 44610.8 -> 44610.734
$ws.Cells.Item(126, 9).Value = 2702  # I126: 2702.1 -> 2702
$ws.Cells.Item(126, 11).Value = 8106  # K126: 8106.299999999999 -> 8106
$ws.Cells.Item(126, 13).Value = -5636  # M126: -5636.299999999999 -> -5636
$ws.Cells.Item(132, 8).Value = 6440.9473  # H132: 7249.5 -> 6440.9473
$ws.Cells.Item(132, 9).Value = 6610.5293  # I132: 7570.9287 -> 6610.5293
$ws.Cells.Item(132, 11).Value = 19831.5879  # K132: 22712.7861 -> 19831.5879
$ws.Cells.Item(132, 13).Value = -17301.5879  # M132: -20182.7861 -> -17301.5879

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 29471.467  # H46: 29508.133 -> 29471.467
$ws.Cells.Item(46, 9).Value = 48163.555  # I46: 61646 -> 48163.555
$ws.Cells.Item(46, 10).Value = 1433.3334  # J46: 1387.5 -> 1433.3334
$ws.Cells.Item(46, 11).Value = 48163.555  # K46: 61646 -> 48163.555
$ws.Cells.Item(46, 12).Value = 1433.3334  # L46: 1387.5 -> 1433.3334
$ws.Cells.Item(46, 13).Value = -47975.555  # M46: -61458 -> -47975.555
$ws.Cells.Item(46, 14).Value = -1809.3334  # N46: -1763.5 -> -1809.3334
$ws.Cells.Item(70, 8).Value = 31900  # H70: 0 -> 31900
$ws.Cells.Item(70, 10).Value = 31900  # J70: 0 -> 31900
$ws.Cells.Item(70, 12).Value = 31900  # L70: 0 -> 31900
$ws.Cells.Item(70, 14).Value = -32440  # N70: None -> -32440
$ws.Cells.Item(73, 8).Value = 31900  # H73: 0 -> 31900
$ws.Cells.Item(73, 10).Value = 31900  # J73: 0 -> 31900
$ws.Cells.Item(73, 12).Value = 31900  # L73: 0 -> 31900
$ws.Cells.Item(73, 14).Value = -33772  # N73: None -> -33772
$ws.Cells.Item(132, 8).Value = 8598.799999999999  # H132: 8598.6 -> 8598.799999999999
$ws.Cells.Item(132, 9).Value = 12499.5  # I132: 19999 -> 12499.5
$ws.Cells.Item(132, 10).Value = 5998.3335  # J132: 5748.5 -> 5998.3335
$ws.Cells.Item(132, 11).Value = 37498.5  # K132: 59997 -> 37498.5
$ws.Cells.Item(132, 12).Value = 17995.0005  # L132: 17245.5 -> 17995.0005
$ws.Cells.Item(132, 13).Value = -34968.5  # M132: -57467 -> -34968.5
$ws.Cells.Item(132, 14).Value = -23055.0005  # N132: -22305.5 -> -23055.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 862970.8  # H132: 4422.2393 -> 862970.8
$ws.Cells.Item(132, 9).Value = 991839  # I132: 4508.175 -> 991839
$ws.Cells.Item(132, 11).Value = 2975517  # K132: 13524.525 -> 2975517
$ws.Cells.Item(132, 13).Value = -2972987  # M132: -10994.525 -> -2972987

Write-Output "Applied price refresh across 8 sheets."